$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.812.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.744.71'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.17%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.601'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.21%  '
$ws.Range("E9").Value = '  -3.16%  '
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("E11").Value = '  -2.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.47'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -18.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.232.16'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.36'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.520.96'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.98%  '
$ws.Range("E16").Value = '  -2.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.749.37'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.13'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '353.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.72'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.99%  '
$ws.Range("E22").Value = '  +0.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.996'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = '  -1.19%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.90%  '
$ws.Range("E28").Value = '  -0.40%  '
$ws.Range("E29").Value = '  -3.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '169.30'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.83%  '
$ws.Range("E32").Value = '  -7.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.10'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.83'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("E36").Value = '  -0.66%  '
$ws.Range("E37").Value = '  -2.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.975'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.11'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.60%  '
$ws.Range("E40").Value = '  -2.78%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '324.95'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.83'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.77%  '
$ws.Range("E44").Value = '  -1.22%  '
$ws.Range("E45").Value = '  -3.36%  '
$ws.Range("E46").Value = '  -1.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '134.80'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.16%  '
$ws.Range("E48").Value = '  -0.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.621'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.32%  '
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.30%  '
